# Auto-generated edit script: apply numeric updates to Moogle_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 12485.143
$ws.Range("I86").Value = 11479.3
$ws.Range("K86").Value = 11479.3
$ws.Range("M86").Value = -10356.3
$ws.Range("H88").Value = 1162.3334
$ws.Range("J88").Value = 994
$ws.Range("L88").Value = 994
$ws.Range("N88").Value = -1806
$ws.Range("H89").Value = 12485.143
$ws.Range("I89").Value = 11479.3
$ws.Range("K89").Value = 57396.5
$ws.Range("M89").Value = -51780.5
$ws.Range("H91").Value = 1162.3334
$ws.Range("J91").Value = 994
$ws.Range("L91").Value = 994
$ws.Range("N91").Value = -3802
$ws.Range("H137").Value = 1832.52
$ws.Range("I137").Value = 1701.5883
$ws.Range("J137").Value = 2110.75
$ws.Range("K137").Value = 5104.7649
$ws.Range("L137").Value = 6332.25
$ws.Range("M137").Value = -2554.7649
$ws.Range("N137").Value = -11432.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 642.2308
$ws.Range("I2").Value = 622.1429000000001
$ws.Range("J2").Value = 726.6
$ws.Range("K2").Value = 622.1429000000001
$ws.Range("L2").Value = 726.6
$ws.Range("M2").Value = -509.1429000000001
$ws.Range("N2").Value = -952.6
$ws.Range("H32").Value = 8972.553
$ws.Range("I32").Value = 5061.0757
$ws.Range("K32").Value = 5061.0757
$ws.Range("M32").Value = -4774.0757
$ws.Range("H33").Value = 16366.333
$ws.Range("I33").Value = 12549.5
$ws.Range("J33").Value = 24000
$ws.Range("K33").Value = 12549.5
$ws.Range("L33").Value = 24000
$ws.Range("M33").Value = -12220.5
$ws.Range("N33").Value = -24658
$ws.Range("H61").Value = 3050.1206
$ws.Range("I61").Value = 2789.6956
$ws.Range("J61").Value = 4048.4167
$ws.Range("K61").Value = 2789.6956
$ws.Range("L61").Value = 4048.4167
$ws.Range("M61").Value = -2577.6956
$ws.Range("N61").Value = -4472.4167
$ws.Range("H88").Value = 4007
$ws.Range("J88").Value = 4007
$ws.Range("L88").Value = 4007
$ws.Range("N88").Value = -4819
$ws.Range("H91").Value = 4007
$ws.Range("J91").Value = 4007
$ws.Range("L91").Value = 4007
$ws.Range("N91").Value = -6815
$ws.Range("H116").Value = 642.2308
$ws.Range("I116").Value = 622.1429000000001
$ws.Range("J116").Value = 726.6
$ws.Range("K116").Value = 622.1429000000001
$ws.Range("L116").Value = 726.6
$ws.Range("M116").Value = 1671.8571
$ws.Range("N116").Value = -5314.6
$ws.Range("H122").Value = 2348.2327
$ws.Range("I122").Value = 1667.6389
$ws.Range("J122").Value = 5848.4287
$ws.Range("K122").Value = 5002.9167
$ws.Range("L122").Value = 17545.2861
$ws.Range("M122").Value = -2552.9167
$ws.Range("N122").Value = -22445.2861
$ws.Range("H136").Value = 3050.1206
$ws.Range("I136").Value = 2789.6956
$ws.Range("J136").Value = 4048.4167
$ws.Range("K136").Value = 8369.086800000001
$ws.Range("L136").Value = 12145.2501
$ws.Range("M136").Value = -5819.086800000001
$ws.Range("N136").Value = -17245.2501

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 642.2308
$ws.Range("I3").Value = 622.1429000000001
$ws.Range("J3").Value = 726.6
$ws.Range("K3").Value = 622.1429000000001
$ws.Range("L3").Value = 726.6
$ws.Range("M3").Value = -508.1429000000001
$ws.Range("N3").Value = -954.6
$ws.Range("H134").Value = 3534.3257
$ws.Range("I134").Value = 2461.6667
$ws.Range("J134").Value = 13992.75
$ws.Range("K134").Value = 7385.000100000001
$ws.Range("L134").Value = 41978.25
$ws.Range("M134").Value = -4850.000100000001
$ws.Range("N134").Value = -47048.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2879.923
$ws.Range("I7").Value = 2912
$ws.Range("K7").Value = 2912
$ws.Range("M7").Value = -2799
$ws.Range("H22").Value = 6735.6
$ws.Range("I22").Value = 7098
$ws.Range("J22").Value = 6192
$ws.Range("K22").Value = 7098
$ws.Range("L22").Value = 6192
$ws.Range("M22").Value = -6748
$ws.Range("N22").Value = -6892
$ws.Range("H31").Value = 11511.637
$ws.Range("I31").Value = 8043.1875
$ws.Range("J31").Value = 14776.059
$ws.Range("K31").Value = 8043.1875
$ws.Range("L31").Value = 14776.059
$ws.Range("M31").Value = -7748.1875
$ws.Range("N31").Value = -15366.059
$ws.Range("H34").Value = 11511.637
$ws.Range("I34").Value = 8043.1875
$ws.Range("J34").Value = 14776.059
$ws.Range("K34").Value = 8043.1875
$ws.Range("L34").Value = 14776.059
$ws.Range("M34").Value = -7841.1875
$ws.Range("N34").Value = -15180.059
$ws.Range("H122").Value = 2322.9
$ws.Range("I122").Value = 2208.5625
$ws.Range("J122").Value = 2453.5715
$ws.Range("K122").Value = 6625.6875
$ws.Range("L122").Value = 7360.7145
$ws.Range("M122").Value = -4175.6875
$ws.Range("N122").Value = -12260.7145
$ws.Range("H132").Value = 5754.121
$ws.Range("I132").Value = 3736.5557
$ws.Range("J132").Value = 14833.167
$ws.Range("K132").Value = 11209.6671
$ws.Range("L132").Value = 44499.501
$ws.Range("M132").Value = -8679.667099999999
$ws.Range("N132").Value = -49559.501
$ws.Range("H134").Value = 3432.7778
$ws.Range("I134").Value = 2872.5
$ws.Range("K134").Value = 8617.5
$ws.Range("M134").Value = -6082.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 424.75
$ws.Range("I29").Value = 400
$ws.Range("K29").Value = 1200
$ws.Range("M29").Value = -923
$ws.Range("H99").Value = 2773.9
$ws.Range("I99").Value = 2773.9
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8321.700000000001
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6075.700000000001
$ws.Range("N99").Value = $null
$ws.Range("H112").Value = 16666.5
$ws.Range("I112").Value = 9999.5
$ws.Range("K112").Value = 29998.5
$ws.Range("M112").Value = -28890.5
$ws.Range("H115").Value = 7499.6665
$ws.Range("I115").Value = 7499.6665
$ws.Range("K115").Value = 22498.9995
$ws.Range("M115").Value = -21323.9995
$ws.Range("H118").Value = 1000
$ws.Range("I118").Value = 1000
$ws.Range("K118").Value = 3000
$ws.Range("M118").Value = -1757
$ws.Range("H120").Value = 16971.375
$ws.Range("I120").Value = 9751.888999999999
$ws.Range("K120").Value = 29255.667
$ws.Range("M120").Value = -24417.667
$ws.Range("H130").Value = 10999
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = $null
$ws.Range("H131").Value = 4547.0625
$ws.Range("J131").Value = 4394.5713
$ws.Range("L131").Value = 13183.7139
$ws.Range("N131").Value = -23263.7139
$ws.Range("H140").Value = 1049.8572
$ws.Range("I140").Value = 1049.8572
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3149.5716
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2030.4284
$ws.Range("N140").Value = $null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5391.609
$ws.Range("I126").Value = 5108.143
$ws.Range("J126").Value = 5832.5557
$ws.Range("K126").Value = 15324.429
$ws.Range("L126").Value = 17497.6671
$ws.Range("M126").Value = -12854.429
$ws.Range("N126").Value = -22437.6671
$ws.Range("H132").Value = 5584.2964
$ws.Range("I132").Value = 4510.48
$ws.Range("K132").Value = 13531.44
$ws.Range("M132").Value = -11001.44

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 170084.5
$ws.Range("I7").Value = 203101.4
$ws.Range("K7").Value = 203101.4
$ws.Range("M7").Value = -202989.4
$ws.Range("H40").Value = 7703.3
$ws.Range("J40").Value = 9187.362999999999
$ws.Range("L40").Value = 9187.362999999999
$ws.Range("N40").Value = -9459.362999999999
$ws.Range("H55").Value = 765.7273
$ws.Range("J55").Value = 1191.5
$ws.Range("L55").Value = 1191.5
$ws.Range("N55").Value = -1537.5
$ws.Range("H58").Value = 10213.714
$ws.Range("I58").Value = 3582.8333
$ws.Range("K58").Value = 3582.8333
$ws.Range("M58").Value = -3322.8333
$ws.Range("H99").Value = 39321
$ws.Range("I99").Value = 29333
$ws.Range("J99").Value = 69285
$ws.Range("K99").Value = 29333
$ws.Range("L99").Value = 69285
$ws.Range("M99").Value = -26338
$ws.Range("N99").Value = -75275
$ws.Range("H100").Value = 4069
$ws.Range("I100").Value = 3083.889
$ws.Range("J100").Value = 8502
$ws.Range("K100").Value = 3083.889
$ws.Range("L100").Value = 8502
$ws.Range("M100").Value = -2542.889
$ws.Range("N100").Value = -9584
$ws.Range("H105").Value = 67000
$ws.Range("J105").Value = 59666.668
$ws.Range("L105").Value = 59666.668
$ws.Range("N105").Value = -66654.66800000001
$ws.Range("H122").Value = 4800.8
$ws.Range("I122").Value = 4001
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 12003
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -9553
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 170084.5
$ws.Range("I126").Value = 203101.4
$ws.Range("K126").Value = 609304.2
$ws.Range("M126").Value = -606834.2
$ws.Range("H132").Value = 3676.6667
$ws.Range("I132").Value = 2955.0625
$ws.Range("J132").Value = 9449.5
$ws.Range("K132").Value = 8865.1875
$ws.Range("L132").Value = 28348.5
$ws.Range("M132").Value = -6335.1875
$ws.Range("N132").Value = -33408.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2616.74
$ws.Range("I132").Value = 2541.4375
$ws.Range("J132").Value = 4424
$ws.Range("K132").Value = 7624.3125
$ws.Range("L132").Value = 13272
$ws.Range("M132").Value = -5094.3125
$ws.Range("N132").Value = -18332

